$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update credential values on rows 2 and 3 ---
# Row 2: swap the Rally login email and password for the new QA creds.
$ws.Range("W2").Value = "RR2tprod@mailinator.com"
$ws.Range("X2").Value = "Testing123!"

# Row 3: swap the Rally login email for the new QA creds and turn it into
# a live mailto hyperlink (matching the styling already used on W2).
$ws.Range("W3").Value = "RR1prodt@mailinator.com"
$ws.Hyperlinks.Add($ws.Range("W3"), "mailto:RR1prodt@mailinator.com")

# Hyperlinks.Add() re-stamps the cell with its own "fresh" hyperlink style
# (no border). Restore the bordered/centered "Hyperlink" look that the rest
# of the sheet's linked cells (e.g. W2) already use.
$ws.Range("W3").Borders.LineStyle = 1
$ws.Range("W3").HorizontalAlignment = -4108
$ws.Range("W3").VerticalAlignment = -4108

# --- Update the active selection (was X2, now S2) ---
$ws.Range("S2").Select() | Out-Null
